# PlayerPerformance_4596.xlsx update:
#  1) Insert a new "Player Info" worksheet as the first sheet, containing
#     the player's basic info (ID, NAME, BATTING_HAND, BOWL_STYLE).
#  2) Keep the existing "ODI Batting" sheet as the second sheet, renaming
#     the MATCH_CARD_LINK column to MATCH_CODE and replacing the full
#     scorecard URLs with just the numeric match code.

$wb = $excel.ActiveWorkbook

# --- existing sheet: ODI Batting -----------------------------------------
$batting = $wb.Worksheets.Item("ODI Batting")

# Rename the MATCH_CARD_LINK header to MATCH_CODE
$batting.Range("D1").Value = "MATCH_CODE"

# Replace each full scorecard URL with just the trailing match code number.
# Prefix with an apostrophe so the numeric-looking text is stored as text,
# matching the rest of the sheet (which stores every value as text).
$batting.Range("D2").Value = "'3946"
$batting.Range("D3").Value = "'3948"
$batting.Range("D4").Value = "'3949"
$batting.Range("D5").Value = "'4698"
$batting.Range("D6").Value = "'4699"
$batting.Range("D7").Value = "'4700"

# --- new sheet: Player Info -----------------------------------------------
# Worksheets.Add() inserts the new sheet right before the currently active
# sheet, which puts it in the first position (matching the target sheet
# order: Player Info, then ODI Batting).
$info = $wb.Worksheets.Add()
$info.Name = "Player Info"

# Match the page margins used throughout the workbook (0.75in sides,
# 1in top/bottom, 0.5in header/footer == 54/72/36 points).
$info.PageSetup.LeftMargin = 54
$info.PageSetup.RightMargin = 54
$info.PageSetup.TopMargin = 72
$info.PageSetup.BottomMargin = 72
$info.PageSetup.HeaderMargin = 36
$info.PageSetup.FooterMargin = 36

$info.Range("A1").Value = "ID"
$info.Range("B1").Value = "NAME"
$info.Range("C1").Value = "BATTING_HAND"
$info.Range("D1").Value = "BOWL_STYLE"

# Match the bold / bordered / centered header styling used on the other
# sheet's header row.
$headerRange = $info.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.Borders.LineStyle = 1
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160

$info.Range("A2").Value = "'4596"
$info.Range("B2").Value = "Ben Matthew Duckett"
$info.Range("C2").Value = "Left Handed"
$info.Range("D2").Value = "Right Arm Off Break"

$info.Range("A1").Select() | Out-Null
